$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 9.008220617589458
$ws.Cells.Item(2, 4).Value = 6.159290980737357
$ws.Cells.Item(2, 5).Value = 12.06062079537978
$ws.Cells.Item(2, 6).Value = 31.27114265021475
$ws.Cells.Item(2, 7).Value = 41.00564469690659
$ws.Cells.Item(2, 8).Value = 17.20488766937506
$ws.Cells.Item(2, 9).Value = 28.1208960237257
$ws.Cells.Item(2, 11).Value = 15.55512219329391
$ws.Cells.Item(2, 12).Value = 9.466769726472325
$ws.Cells.Item(2, 13).Value = 17.58673166196473
$ws.Cells.Item(2, 14).Value = 18.65459400662565

$ws.Cells.Item(3, 3).Value = 8.99029258589918
$ws.Cells.Item(3, 4).Value = 6.160487358684868
$ws.Cells.Item(3, 5).Value = 12.07776283180736
$ws.Cells.Item(3, 6).Value = 31.25138441817989
$ws.Cells.Item(3, 7).Value = 40.94685972615206
$ws.Cells.Item(3, 8).Value = 17.24851011020611
$ws.Cells.Item(3, 9).Value = 28.1548618585667
$ws.Cells.Item(3, 11).Value = 15.18128324418266
$ws.Cells.Item(3, 12).Value = 9.48993581352199
$ws.Cells.Item(3, 13).Value = 17.44684810814705
$ws.Cells.Item(3, 14).Value = 18.71873789269169

$ws.Cells.Item(4, 3).Value = 8.980714304919706
$ws.Cells.Item(4, 4).Value = 6.16139182307613
$ws.Cells.Item(4, 5).Value = 12.09016701775963
$ws.Cells.Item(4, 6).Value = 31.24860317656063
$ws.Cells.Item(4, 7).Value = 40.92498149230146
$ws.Cells.Item(4, 8).Value = 17.27887976525138
$ws.Cells.Item(4, 9).Value = 28.18289115940418
$ws.Cells.Item(4, 11).Value = 14.94994757091527
$ws.Cells.Item(4, 12).Value = 9.505222096788685
$ws.Cells.Item(4, 13).Value = 17.36363831509599
$ws.Cells.Item(4, 14).Value = 18.76002650624305

$ws.Cells.Item(5, 3).Value = 8.97717330360055
$ws.Cells.Item(5, 4).Value = 6.161803398900703
$ws.Cells.Item(5, 5).Value = 12.09569422387516
$ws.Cells.Item(5, 6).Value = 31.24981955511467
$ws.Cells.Item(5, 7).Value = 40.91963991911581
$ws.Cells.Item(5, 8).Value = 17.29215497775573
$ws.Cells.Item(5, 9).Value = 28.19611259890121
$ws.Cells.Item(5, 11).Value = 14.85536547216185
$ws.Cells.Item(5, 12).Value = 9.511718769589049
$ws.Cells.Item(5, 13).Value = 17.33043382301307
$ws.Cells.Item(5, 14).Value = 18.77733217196763

$ws.Cells.Item(6, 3).Value = 8.976607275556962
$ws.Cells.Item(6, 4).Value = 6.161874346067743
$ws.Cells.Item(6, 5).Value = 12.0966405406589
$ws.Cells.Item(6, 6).Value = 31.25016337100468
$ws.Cells.Item(6, 7).Value = 40.9189687141992
$ws.Cells.Item(6, 8).Value = 17.29441356687443
$ws.Cells.Item(6, 9).Value = 28.19841652754202
$ws.Cells.Item(6, 11).Value = 14.83964563434083
$ws.Cells.Item(6, 12).Value = 9.51281369521972
$ws.Cells.Item(6, 13).Value = 17.32496363837388
$ws.Cells.Item(6, 14).Value = 18.78023481127968

$ws.Cells.Item(7, 3).Value = 8.980665079663034
$ws.Cells.Item(7, 4).Value = 6.16139719927151
$ws.Cells.Item(7, 5).Value = 12.09023964700674
$ws.Cells.Item(7, 6).Value = 31.24861007050551
$ws.Cells.Item(7, 7).Value = 40.92489498713071
$ws.Cells.Item(7, 8).Value = 17.27905516111425
$ws.Cells.Item(7, 9).Value = 28.18306219004309
$ws.Cells.Item(7, 11).Value = 14.94867307231326
$ws.Cells.Item(7, 12).Value = 9.505308630166319
$ws.Cells.Item(7, 13).Value = 17.36318761723867
$ws.Cells.Item(7, 14).Value = 18.76025795021085

$ws.Cells.Item(8, 3).Value = 9.001744037492255
$ws.Cells.Item(8, 4).Value = 6.159668399633317
$ws.Cells.Item(8, 5).Value = 12.06614138181325
$ws.Cells.Item(8, 6).Value = 31.26238904461347
$ws.Cells.Item(8, 7).Value = 40.98242587863619
$ws.Cells.Item(8, 8).Value = 17.21918325546438
$ws.Cells.Item(8, 9).Value = 28.13111583061048
$ws.Cells.Item(8, 11).Value = 15.42667827092358
$ws.Cells.Item(8, 12).Value = 9.474537103221154
$ws.Cells.Item(8, 13).Value = 17.53796235788696
$ws.Cells.Item(8, 14).Value = 18.67631651062623

$ws.Cells.Item(9, 3).Value = 9.054292391097814
$ws.Cells.Item(9, 4).Value = 6.157614147686255
$ws.Cells.Item(9, 5).Value = 12.03379428642691
$ws.Cells.Item(9, 6).Value = 31.36357879296433
$ws.Cells.Item(9, 7).Value = 41.2078901939177
$ws.Cells.Item(9, 8).Value = 17.13031779199935
$ws.Cells.Item(9, 9).Value = 28.08635283959902
$ws.Cells.Item(9, 11).Value = 16.34390441692847
$ws.Cells.Item(9, 12).Value = 9.422609661213059
$ws.Cells.Item(9, 13).Value = 17.90056297165981
$ws.Cells.Item(9, 14).Value = 18.52674884227107

$ws.Cells.Item(10, 3).Value = 9.099545313946786
$ws.Cells.Item(10, 4).Value = 6.156903566805378
$ws.Cells.Item(10, 5).Value = 12.01911720662644
$ws.Cells.Item(10, 6).Value = 31.4829811014619
$ws.Cells.Item(10, 7).Value = 41.4417410349966
$ws.Cells.Item(10, 8).Value = 17.08255355629827
$ws.Cells.Item(10, 9).Value = 28.08847709253378
$ws.Cells.Item(10, 11).Value = 16.99813116306848
$ws.Cells.Item(10, 12).Value = 9.389571456019967
$ws.Cells.Item(10, 13).Value = 18.17715075465819
$ws.Cells.Item(10, 14).Value = 18.4259377705385

$ws.Cells.Item(11, 3).Value = 9.12152948566836
$ws.Cells.Item(11, 4).Value = 6.156750458249089
$ws.Cells.Item(11, 5).Value = 12.01441193498113
$ws.Cells.Item(11, 6).Value = 31.54700686633898
$ws.Cells.Item(11, 7).Value = 41.56274825699441
$ws.Cells.Item(11, 8).Value = 17.06465381768086
$ws.Cells.Item(11, 9).Value = 28.09707122181425
$ws.Cells.Item(11, 11).Value = 17.29008035576229
$ws.Cells.Item(11, 12).Value = 9.375648125744842
$ws.Cells.Item(11, 13).Value = 18.30476500384343
$ws.Cells.Item(11, 14).Value = 18.3820273285446

$ws.Cells.Item(12, 3).Value = 9.130051136097642
$ws.Cells.Item(12, 4).Value = 6.156716673404091
$ws.Cells.Item(12, 5).Value = 12.01291333927541
$ws.Cells.Item(12, 6).Value = 31.57263785528859
$ws.Cells.Item(12, 7).Value = 41.61065160476003
$ws.Cells.Item(12, 8).Value = 17.05842789521701
$ws.Cells.Item(12, 9).Value = 28.10142295894459
$ws.Cells.Item(12, 11).Value = 17.39971365864286
$ws.Cells.Item(12, 12).Value = 9.370534487948012
$ws.Cells.Item(12, 13).Value = 18.35331076196786
$ws.Cells.Item(12, 14).Value = 18.36567840413845

$ws.Cells.Item(13, 3).Value = 9.128207168912651
$ws.Cells.Item(13, 4).Value = 6.15672287783039
$ws.Cells.Item(13, 5).Value = 12.01322349948549
$ws.Cells.Item(13, 6).Value = 31.56705631533628
$ws.Cells.Item(13, 7).Value = 41.60024263551021
$ws.Cells.Item(13, 8).Value = 17.05974416758331
$ws.Cells.Item(13, 9).Value = 28.10043693315343
$ws.Cells.Item(13, 11).Value = 17.37614477029491
$ws.Cells.Item(13, 12).Value = 9.371628740756851
$ws.Cells.Item(13, 13).Value = 18.34284632264211
$ws.Cells.Item(13, 14).Value = 18.36918704845423

$ws.Cells.Item(14, 3).Value = 9.122226650203938
$ws.Cells.Item(14, 4).Value = 6.156747195412315
$ws.Cells.Item(14, 5).Value = 12.01428297111449
$ws.Cells.Item(14, 6).Value = 31.54908783597205
$ws.Cells.Item(14, 7).Value = 41.56664774225224
$ws.Cells.Item(14, 8).Value = 17.06413052827888
$ws.Cells.Item(14, 9).Value = 28.09740725045894
$ws.Cells.Item(14, 11).Value = 17.2991190073033
$ws.Cells.Item(14, 12).Value = 9.375224241747308
$ws.Cells.Item(14, 13).Value = 18.30875467463208
$ws.Cells.Item(14, 14).Value = 18.3806767081229

$ws.Cells.Item(15, 3).Value = 9.118588895289575
$ws.Cells.Item(15, 4).Value = 6.156765233393811
$ws.Cells.Item(15, 5).Value = 12.01496879722163
$ws.Cells.Item(15, 6).Value = 31.53826175575578
$ws.Cells.Item(15, 7).Value = 41.54634010113229
$ws.Cells.Item(15, 8).Value = 17.06688927970852
$ws.Cells.Item(15, 9).Value = 28.09569438444302
$ws.Cells.Item(15, 11).Value = 17.2518154513562
$ws.Cells.Item(15, 12).Value = 9.377447267660601
$ws.Cells.Item(15, 13).Value = 18.28790024250212
$ws.Cells.Item(15, 14).Value = 18.38775075433209

$ws.Cells.Item(16, 3).Value = 9.098136367993481
$ws.Cells.Item(16, 4).Value = 6.156916973143751
$ws.Cells.Item(16, 5).Value = 12.01946435357595
$ws.Cells.Item(16, 6).Value = 31.4789914036003
$ws.Cells.Item(16, 7).Value = 41.43412540592713
$ws.Cells.Item(16, 8).Value = 17.08380054150473
$ws.Cells.Item(16, 9).Value = 28.08806899715176
$ws.Cells.Item(16, 11).Value = 16.97892829741937
$ws.Cells.Item(16, 12).Value = 9.390503614707168
$ws.Cells.Item(16, 13).Value = 18.16884380020094
$ws.Cells.Item(16, 14).Value = 18.42884653230504

$ws.Cells.Item(17, 3).Value = 9.085944534923994
$ws.Cells.Item(17, 4).Value = 6.157053442940124
$ws.Cells.Item(17, 5).Value = 12.02272697325299
$ws.Cells.Item(17, 6).Value = 31.44511100957267
$ws.Cells.Item(17, 7).Value = 41.36901675459386
$ws.Cells.Item(17, 8).Value = 17.09515693783029
$ws.Cells.Item(17, 9).Value = 28.08534535131885
$ws.Cells.Item(17, 11).Value = 16.80999168693026
$ws.Cells.Item(17, 12).Value = 9.398796365708689
$ws.Cells.Item(17, 13).Value = 18.09623963614782
$ws.Cells.Item(17, 14).Value = 18.45455574904695

$ws.Cells.Item(18, 3).Value = 9.079064014492094
$ws.Cells.Item(18, 4).Value = 6.157147966991226
$ws.Cells.Item(18, 5).Value = 12.02478912728929
$ws.Cells.Item(18, 6).Value = 31.42653851912881
$ws.Cells.Item(18, 7).Value = 41.33294639883746
$ws.Cells.Item(18, 8).Value = 17.10204909118735
$ws.Cells.Item(18, 9).Value = 28.08449677739632
$ws.Cells.Item(18, 11).Value = 16.71229678394425
$ws.Cells.Item(18, 12).Value = 9.40367024264339
$ws.Cells.Item(18, 13).Value = 18.05465063554491
$ws.Cells.Item(18, 14).Value = 18.46952653546644

$ws.Cells.Item(19, 3).Value = 9.076757171523347
$ws.Cells.Item(19, 4).Value = 6.157182732755665
$ws.Cells.Item(19, 5).Value = 12.02551921859007
$ws.Cells.Item(19, 6).Value = 31.42040755923029
$ws.Cells.Item(19, 7).Value = 41.3209709710239
$ws.Cells.Item(19, 8).Value = 17.10444447274854
$ws.Cells.Item(19, 9).Value = 28.08433275209974
$ws.Cells.Item(19, 11).Value = 16.67913185872897
$ws.Cells.Item(19, 12).Value = 9.405338340200737
$ws.Cells.Item(19, 13).Value = 18.04059980220459
$ws.Cells.Item(19, 14).Value = 18.47462694584752

$ws.Cells.Item(20, 3).Value = 9.087228756154973
$ws.Cells.Item(20, 4).Value = 6.157037258390605
$ws.Cells.Item(20, 5).Value = 12.02236045669735
$ws.Cells.Item(20, 6).Value = 31.44862304405392
$ws.Cells.Item(20, 7).Value = 41.37580517319442
$ws.Cells.Item(20, 8).Value = 17.09391073152079
$ws.Cells.Item(20, 9).Value = 28.08556097194178
$ws.Cells.Item(20, 11).Value = 16.82803062771164
$ws.Cells.Item(20, 12).Value = 9.397902815124397
$ws.Cells.Item(20, 13).Value = 18.10395104234788
$ws.Cells.Item(20, 14).Value = 18.45179997405843

$ws.Cells.Item(21, 3).Value = 9.123977970480151
$ws.Cells.Item(21, 4).Value = 6.156739398279663
$ws.Cells.Item(21, 5).Value = 12.01396409577732
$ws.Cells.Item(21, 6).Value = 31.554328095909
$ws.Cells.Item(21, 7).Value = 41.57645911345595
$ws.Cells.Item(21, 8).Value = 17.06282714451202
$ws.Cells.Item(21, 9).Value = 28.09826736151641
$ws.Cells.Item(21, 11).Value = 17.32176916478427
$ws.Cells.Item(21, 12).Value = 9.374163847086818
$ws.Cells.Item(21, 13).Value = 18.31876251405881
$ws.Cells.Item(21, 14).Value = 18.37729435433173

$ws.Cells.Item(22, 3).Value = 9.149140165974796
$ws.Cells.Item(22, 4).Value = 6.15668562676063
$ws.Cells.Item(22, 5).Value = 12.01012702596207
$ws.Cells.Item(22, 6).Value = 31.63148381337079
$ws.Cells.Item(22, 7).Value = 41.71971164306868
$ws.Cells.Item(22, 8).Value = 17.04573228456306
$ws.Cells.Item(22, 9).Value = 28.11296733713653
$ws.Cells.Item(22, 11).Value = 17.6390406209645
$ws.Cells.Item(22, 12).Value = 9.359574722003758
$ws.Cells.Item(22, 13).Value = 18.46042531091258
$ws.Cells.Item(22, 14).Value = 18.33022632104907

$ws.Cells.Item(23, 3).Value = 9.135607411674991
$ws.Cells.Item(23, 4).Value = 6.156701522586971
$ws.Cells.Item(23, 5).Value = 12.01202404886574
$ws.Cells.Item(23, 6).Value = 31.58956978104221
$ws.Cells.Item(23, 7).Value = 41.64215525447856
$ws.Cells.Item(23, 8).Value = 17.05456094674599
$ws.Cells.Item(23, 9).Value = 28.10453657331474
$ws.Cells.Item(23, 11).Value = 17.47023566942666
$ws.Cells.Item(23, 12).Value = 9.367276580850719
$ws.Cells.Item(23, 13).Value = 18.38471310048316
$ws.Cells.Item(23, 14).Value = 18.3551990758755

$ws.Cells.Item(24, 3).Value = 9.086647758454538
$ws.Cells.Item(24, 4).Value = 6.157044525334268
$ws.Cells.Item(24, 5).Value = 12.02252557809972
$ws.Cells.Item(24, 6).Value = 31.44703243090706
$ws.Cells.Item(24, 7).Value = 41.37273188643817
$ws.Cells.Item(24, 8).Value = 17.09447301008907
$ws.Cells.Item(24, 9).Value = 28.08546125552346
$ws.Cells.Item(24, 11).Value = 16.81987699282718
$ws.Cells.Item(24, 12).Value = 9.398306458318311
$ws.Cells.Item(24, 13).Value = 18.10046423899538
$ws.Cells.Item(24, 14).Value = 18.4530452672816

$ws.Cells.Item(25, 3).Value = 9.038894861639788
$ws.Cells.Item(25, 4).Value = 6.158028494200426
$ws.Cells.Item(25, 5).Value = 12.04094863478784
$ws.Cells.Item(25, 6).Value = 31.32827182106893
$ws.Cells.Item(25, 7).Value = 41.13486855339437
$ws.Cells.Item(25, 8).Value = 17.15128911576582
$ws.Cells.Item(25, 9).Value = 28.09232510061727
$ws.Cells.Item(25, 11).Value = 16.09871824198532
$ws.Cells.Item(25, 12).Value = 9.435758292973386
$ws.Cells.Item(25, 13).Value = 17.80053874851243
$ws.Cells.Item(25, 14).Value = 18.5656103053733
